$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resources")

# Row 2 - population
$ws.Range("E2").Value = 12000
$ws.Range("F2").Value = 18000

# Row 3 - metalElements
$ws.Range("E3").Value = 1200
$ws.Range("F3").Value = 1800

# Row 4 - timber
$ws.Range("E4").Value = 8250
$ws.Range("F4").Value = 27000

# Row 5 - landArea
$ws.Range("E5").Value = 15000
$ws.Range("F5").Value = 30000

# Row 6 - water
$ws.Range("B6").Value = 2
$ws.Range("E6").Value = 40000
$ws.Range("F6").Value = 55000

# Row 7 - metalAlloys
$ws.Range("E7").Value = 700
$ws.Range("F7").Value = 1300

# Row 8 - electronics
$ws.Range("B8").Value = 4300
$ws.Range("E8").Value = 700
$ws.Range("F8").Value = 1300

# Row 9 - housing
$ws.Range("E9").Value = 14000
$ws.Range("F9").Value = 20000

# Row 10 - food
$ws.Range("E10").Value = 14000
$ws.Range("F10").Value = 20000

# Row 11 - metalAlloysWaste
$ws.Range("E11").Value = 1500

# Row 12 - housingWaste
$ws.Range("E12").Value = 1500

# Row 13 - electronicsWaste
$ws.Range("E13").Value = 1500

# Row 14 - foodWaste
$ws.Range("E14").Value = 1500
